## Rename the header row strings so each column header carries the
## corresponding AHB format-version suffix instead of the old "_old"/"_new"
## markers, then wrap the data range in a table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells: "_old" -> "_FV2310", "_new" -> "_FV2404" ------
$lastCol = 21
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value2
    if ($v -ne $null) {
        $newValue = $v -replace '_old$', '_FV2310'
        $newValue = $newValue -replace '_new$', '_FV2404'
        $cell.Value2 = $newValue
    }
}

# --- 2. Turn the used range into an Excel Table ("Table1") ----------------
$lastRow = $ws.UsedRange.Rows.Count
$tableRange = $ws.Range("A1").Resize($lastRow, $lastCol)
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row (split below row 1) --------------------------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
